# Update (Analyze PO & Forecast)
$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Update MyForecast (column D) values on "Forecast Comparison" sheet
$wsForecast.Range("D8").Value  = 458
$wsForecast.Range("D9").Value  = 335
$wsForecast.Range("D10").Value = 303
$wsForecast.Range("D11").Value = 293
$wsForecast.Range("D12").Value = 514
$wsForecast.Range("D14").Value = 500
$wsForecast.Range("D15").Value = 345
$wsForecast.Range("D16").Value = 450
$wsForecast.Range("D17").Value = 471

# Update summary metrics on "Summary" sheet.
# These cells hold numeric-looking / date-looking text, so force text
# formatting first to keep them stored as text (matching the source file)
# rather than letting Excel auto-convert them to numbers/dates.
$wsSummary.Range("B9").NumberFormat  = "@"
$wsSummary.Range("B9").Value         = "7829"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value        = "4450"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value        = "293"

$wsSummary.Range("B15").NumberFormat = "@"
$wsSummary.Range("B15").Value        = "2025-03-30"
